$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Avoid recalculating dependent formulas (e.g. column B's SUM) so their
# cached values stay as they were before this edit, matching Excel's
# "manual calculation" style in-place cell edit.
$excel.Calculation = -4135

# Update the header for column B (total chai sales) to the new Korean text
$ws.Range("B1").Value = "총 차이 판매액(단위)"

# D7 changes from a numeric value (548) to a text value "5:48"
$ws.Range("D7").Value = "5:48"
